$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TodasLasLineas")

# New rows of data to append (Linea / Caso names + numeric values)
$newData = @(
    @("ABM_Usuarios_MIX", 1162816939),
    @("ABM_Usuarios_POS", 1145642605),
    @("ABM_Usuarios_PRE", 1162676705),
    @("Actualizacion_de_Datos_MIX", 1162816939),
    @("Actualizacion_de_Datos_POS", 1145642605),
    @("Actualizacion_de_Datos_PRE", 1162676705)
)

$startRow = 30
$endRow = $startRow + $newData.Count - 1

for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newData[$i][0]
    $ws.Cells.Item($row, 2).Value = $newData[$i][1]
}

# Carry the number/style formatting from the last existing data row (B29)
# down onto the newly appended B30:B35 cells, matching the workbook's
# existing "s=4" numeric style used throughout column B.
$ws.Range("B29").Copy()
$ws.Range("B" + $startRow + ":B" + $endRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update view state to match the final saved view
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("C35").Select()
